$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4, pushing existing rows 4-28 down to 5-29.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the "Vacuum chamber / Mother volume radius" entry.
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = "Facility"
$ws.Range("C4").Value = "Global"
$ws.Range("D4").Value = "Vacuum chamber"
$ws.Range("E4").Value = "Mother volume radius"
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = "m"

# Match style used elsewhere in this "Facility" block (same as row 3 prior to insert).
$ws.Range("A4:H4").Style = $ws.Range("A3:H3").Style

# Update selection to match the recorded cursor position after the edit.
$ws.Range("F4").Select()
